# Updates the "cryptos" price/volume list to refreshed values.
# All Price (D) and Volume(1h) (E) cells are text (e.g. "42.497.06", "  -1.47%  "),
# so NumberFormat is forced to "@" (Text) before assignment to prevent Excel
# from auto-converting numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.497.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.186.43'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.18%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.64'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.97'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.26'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.45%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.515.45'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.16'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.190.44'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.770'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.406.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.44%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.76'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.74'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.47'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -9.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.27%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.44'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.64%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.36'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.00'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.71%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.14'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.41%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.21'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.33%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.53%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.58'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.16'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.75%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '58.56'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.35'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0971'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.10%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.17%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.12'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.38%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.71%  '

# Row 29: PancakeSwap -> InjectiveProtocol (with updated price/volume)
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.48'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.42%  '

# Row 30: InjectiveProtocol -> PancakeSwap (with updated price/volume)
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.74%  '
